# Daily refresh of the cryptos price/volume table (GitHub Actions job).
# Values below mirror the new coinranking.com snapshot:
#  - D/E columns are updated to the latest Price / Volume(1h) readings.
#  - A couple of price columns ("D5", "D7", ...) look like plain numbers
#    (single decimal point) so a leading apostrophe is used to force them
#    to stay text, matching the original inline-string cell type instead
#    of letting Excel auto-convert them to numbers.
#  - Rows 34/35 and 43/44 had their rankings swap places, so B/C/D/E are
#    rewritten for those rows too.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '37.395.19'
$ws.Range('E2').Value = '  +4.97%  '
$ws.Range('D3').Value = '2.052.47'
$ws.Range('E3').Value = '  +3.47%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '''252.39'
$ws.Range('E5').Value = '  +2.90%  '
$ws.Range('E6').Value = '  +2.37%  '
$ws.Range('D7').Value = '''65.36'
$ws.Range('E7').Value = '  +13.66%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  +6.26%  '
$ws.Range('D10').Value = '''59.15'
$ws.Range('E10').Value = '  +1.27%  '
$ws.Range('D11').Value = '''0.0768'
$ws.Range('E11').Value = '  +4.50%  '
$ws.Range('E12').Value = '  +1.37%  '
$ws.Range('D13').Value = '''0.926'
$ws.Range('E13').Value = '  -2.13%  '
$ws.Range('D14').Value = '''14.83'
$ws.Range('E14').Value = '  +2.55%  '
$ws.Range('D15').Value = '''22.44'
$ws.Range('E15').Value = '  +25.29%  '
$ws.Range('D16').Value = '2.350.55'
$ws.Range('E16').Value = '  +3.54%  '
$ws.Range('D17').Value = '''5.57'
$ws.Range('E17').Value = '  +5.39%  '
$ws.Range('D18').Value = '2.052.43'
$ws.Range('E18').Value = '  +3.73%  '
$ws.Range('D19').Value = '37.308.99'
$ws.Range('E19').Value = '  +4.89%  '
$ws.Range('D20').Value = '''73.46'
$ws.Range('E20').Value = '  +3.02%  '
$ws.Range('E21').Value = '  +3.73%  '
$ws.Range('D22').Value = '''5.48'
$ws.Range('E22').Value = '  +6.21%  '
$ws.Range('D23').Value = '''239.75'
$ws.Range('E23').Value = '  +2.90%  '
$ws.Range('E24').Value = '  +4.41%  '
$ws.Range('E25').Value = '  -0.14%  '
$ws.Range('E26').Value = '  +4.97%  '
$ws.Range('D27').Value = '''10.06'
$ws.Range('E27').Value = '  +10.17%  '
$ws.Range('D28').Value = '''161.85'
$ws.Range('E28').Value = '  -1.65%  '
$ws.Range('E29').Value = '  +3.99%  '
$ws.Range('D30').Value = '''0.122'
$ws.Range('E30').Value = '  +28.48%  '
$ws.Range('E31').Value = '  +8.26%  '
$ws.Range('E32').Value = '  +2.69%  '
$ws.Range('D33').Value = '''1.21'
$ws.Range('E33').Value = '  +8.81%  '
$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D34').Value = '''4.68'
$ws.Range('E34').Value = '  +8.22%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').Value = '''0.0626'
$ws.Range('E35').Value = '  +5.75%  '
$ws.Range('E36').Value = '  +1.22%  '
$ws.Range('E37').Value = '  +0.03%  '
$ws.Range('D38').Value = '''1.84'
$ws.Range('E38').Value = '  +4.30%  '
$ws.Range('E39').Value = '  +14.94%  '
$ws.Range('D40').Value = '''3.04'
$ws.Range('E40').Value = '  +35.56%  '
$ws.Range('E41').Value = '  +17.06%  '
$ws.Range('E42').Value = '  +4.54%  '
$ws.Range('B43').Value = 'InjectiveProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D43').Value = '''17.74'
$ws.Range('E43').Value = '  +10.57%  '
$ws.Range('B44').Value = 'TrustWalletToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D44').Value = '''1.24'
$ws.Range('E44').Value = '  +1.68%  '
$ws.Range('E45').Value = '  +5.67%  '
$ws.Range('D46').Value = '''0.0218'
$ws.Range('E46').Value = '  +3.22%  '
$ws.Range('D47').Value = '''96.72'
$ws.Range('E47').Value = '  +5.24%  '
$ws.Range('D48').Value = '''7.89'
$ws.Range('E48').Value = '  +2.60%  '
$ws.Range('D49').Value = '1.417.69'
$ws.Range('E49').Value = '  +3.50%  '
$ws.Range('E50').Value = '  +1.88%  '
$ws.Range('D51').Value = '''46.56'
$ws.Range('E51').Value = '  -0.72%  '
